# Add a new data row (row 6) to the HZNP "Bag" sentiment scoring sheet,
# matching the format of the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (e.g. the date number format) from the row above down
# into the new row so we don't introduce a brand new cell style.
$ws.Range("A5:N5").Copy()
$ws.Range("A6:N6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new row's values.
$ws.Range("A6").Value = 42608.890162037038
$ws.Range("B6").Value = -4
$ws.Range("C6").Value = 51
$ws.Range("D6").Value = 47
$ws.Range("E6").Value = 35
$ws.Range("F6").Value = 64
$ws.Range("G6").Value = 17507
$ws.Range("H6").Value = 24604
$ws.Range("I6").Value = 2819
$ws.Range("J6").Value = 235
$ws.Range("K6").Value = 216
$ws.Range("L6").Value = 6
$ws.Range("M6").Value = 11
$ws.Range("N6").Value = "Bag"
